# Adds the "Do pre-market engagement" / "Name your project" test data columns
# (M: DoPreMarketEngagement_PageTitle, N: DoPreMarketEngagementPageContentButtonStatus)
# to the OneFCFlowTestData sheet.
#
# Cell values are written in this specific order so the workbook's shared
# string table is built up in the same sequence as the source edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$statusText = "1. Project and colleagues!Name your project!OPTIONAL!Change who is going to lead the procurement!OPTIONAL!Add colleagues to your project!OPTIONAL!2. Build your RfI!Build your RfI!TO DO!3. Review and publish your RfI!Upload documents!CANNOT START YET!See the suppliers who will receive your RfI!CANNOT START YET!Your RfI timeline!CANNOT START YET!Review and publish your RfI!CANNOT START YET"

$ws.Range("N1").Value = "DoPreMarketEngagementPageContentButtonStatus"
$ws.Range("M2").Value = "Do pre-market engagement"
$ws.Range("N2").Value = $statusText
$ws.Range("M1").Value = "DoPreMarketEngagement_PageTitle"
$ws.Range("M3").Value = "Do pre-market engagement"
$ws.Range("N3").Value = $statusText

# --- Match the "vertical top + wrap text" formatting used throughout columns J/L ---
$ws.Range("M1:N3").WrapText = $true
$ws.Range("M1:N3").VerticalAlignment = -4160

# --- Column widths, as close as this engine's width quantization allows to
#     the bestFit widths (30.453125 / 44.1796875) used by the source file ---
$ws.Columns("M").ColumnWidth = 29.59
$ws.Columns("N").ColumnWidth = 43.25

# --- Update selection / active cell to N2, matching the saved view ---
$null = $ws.Range("N2").Select()
